$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string but must remain text
# (matches source data which stores these as inline text, not numbers).
# Force Text number format first so Excel does not coerce them to doubles.
$textCells = @("D4","D5","D6","D7","D8","D9","D11","D12","D13","D15","D17","D18","D19","D20","D21","D22","D24","D25","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.953.05'
$ws.Range('E2').Value = '  -2.12%  '
$ws.Range('D3').Value = '1.866.91'
$ws.Range('E3').Value = '  -2.82%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '312.33'
$ws.Range('E5').Value = '  -1.13%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').Value = '0.4983'
$ws.Range('E7').Value = '  -2.88%  '
$ws.Range('D8').Value = '0.3822'
$ws.Range('E8').Value = '  -4.22%  '
$ws.Range('D9').Value = '0.08929'
$ws.Range('E9').Value = '  -8.74%  '
$ws.Range('E10').Value = '  -2.65%  '
$ws.Range('D11').Value = '41.50'
$ws.Range('E11').Value = '  -1.43%  '
$ws.Range('D12').Value = '6.375'
$ws.Range('E12').Value = '  -1.86%  '
$ws.Range('D13').Value = '20.68'
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('D14').Value = '1.864.26'
$ws.Range('E14').Value = '  -2.41%  '
$ws.Range('D15').Value = '7.225'
$ws.Range('E15').Value = '  -2.80%  '
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').Value = '0.00001099'
$ws.Range('E17').Value = '  -3.01%  '
$ws.Range('D18').Value = '91.04'
$ws.Range('E18').Value = '  -3.54%  '
$ws.Range('D19').Value = '0.06662'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').Value = '18.07'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').Value = '6.120'
$ws.Range('E22').Value = '  -2.81%  '
$ws.Range('D23').Value = '27.978.09'
$ws.Range('E23').Value = '  -2.34%  '
$ws.Range('D24').Value = '11.50'
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('D25').Value = '2.282'
$ws.Range('E25').Value = '  -0.91%  '
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').Value = '2.077.18'
$ws.Range('E26').Value = '  -2.81%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.494'
$ws.Range('E27').Value = '  -7.99%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '158.12'
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '20.70'
$ws.Range('E29').Value = '  -2.64%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '125.97'
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.1059'
$ws.Range('E31').Value = '  -1.20%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '1.056'
$ws.Range('E32').Value = '  -5.41%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '5.573'
$ws.Range('E33').Value = '  -2.44%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '3.582'
$ws.Range('E34').Value = '  -1.25%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').Value = '9.386'
$ws.Range('E35').Value = '  -4.37%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '0.06552'
$ws.Range('E36').Value = '  -2.78%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.02410'
$ws.Range('E37').Value = '  -1.55%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = '0.2189'
$ws.Range('E38').Value = '  -1.89%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '1.289'
$ws.Range('E39').Value = '  +8.30%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '1.195'
$ws.Range('E40').Value = '  -5.93%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').Value = '11.58'
$ws.Range('E41').Value = '  -1.18%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.6363'
$ws.Range('E42').Value = '  -1.54%  '
$ws.Range('B43').Value = 'InternetComputer(DFINITY)'
$ws.Range('C43').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D43').Value = '4.892'
$ws.Range('E43').Value = '  -3.61%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = '1.002'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.6004'
$ws.Range('E45').Value = '  -1.25%  '
$ws.Range('D46').Value = '13.13'
$ws.Range('E46').Value = '  -4.25%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').Value = '3.676'
$ws.Range('E47').Value = '  -2.78%  '
$ws.Range('B48').Value = 'WEMIXTOKEN'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '1.275'
$ws.Range('E48').Value = '  -0.85%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '1.993'
$ws.Range('E49').Value = '  -3.82%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').Value = '1.219'
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = '121.03'
$ws.Range('E51').Value = '  -2.47%  '
